$wb = $excel.ActiveWorkbook

# The data/formatting changes target the "harsha" worksheet (3rd sheet).
$ws = $wb.Worksheets.Item("harsha")

# Make "harsha" the active sheet (drives workbookView activeTab + this
# sheet's sheetView tabSelected; "karthik" loses tabSelected automatically).
$ws.Activate()

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "role"
$ws.Range("F1").Value = "skills"

# Row 2 - harsha's own record
$ws.Range("A2").Value = "harsha1222222"
$ws.Range("B2").Value = "harsha12344690"
$ws.Range("C2").Value = "harsha12389@yopmail.com"
$ws.Range("D2").Value = 12345678
$ws.Range("E2").Value = "HR"

# Row 3 - vardhan's record
$ws.Range("A3").Value = "vardhan123444"
$ws.Range("B3").Value = "vardhan12356666"
$ws.Range("C3").Value = "vardhan124444@yopmail.com"
$ws.Range("D3").Value = 12345678
$ws.Range("E3").Value = "Subject Expert"
$ws.Range("F3").Value = "html"
$ws.Range("G3").Value = "css"

# Column widths for A:D (closest values achievable given the runtime's
# character-width rounding; targets ~15.63, 17.73, 28.09, 15.54 chars)
$ws.Columns.Item(1).ColumnWidth = 14.833333333333332
$ws.Columns.Item(2).ColumnWidth = 16.833333333333336
$ws.Columns.Item(3).ColumnWidth = 27.333333333333336
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666

# Leave the selection on E5, as in the target workbook.
$ws.Range("E5").Select()
